$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.003676399999999996
$ws.Range("E2").Value = 0.3723754890751533
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3669021
$ws.Range("L2").Value = 0.5961429402307628
$ws.Range("M2").Value = 0.08239116666666667
$ws.Range("N2").Value = 12.8635987649262
$ws.Range("O2").Value = 3.475337169561451

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.07324717743582368
$ws.Range("E2").Value = 0.3629913768432933
$ws.Range("I2").Value = 0.57664568317287
$ws.Range("L2").Value = 0.2773383432630344
$ws.Range("M2").Value = 0.08143291666666665
$ws.Range("N2").Value = 9.162835765309527
$ws.Range("O2").Value = 3.84946491064901

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.09755888307768322
$ws.Range("B2").Value = 0.0231476116393296
$ws.Range("E2").Value = 0.157117409245859
$ws.Range("I2").Value = 0.5170865166666656
$ws.Range("M2").Value = 0.04734183333333336
$ws.Range("N2").Value = 8.874638901356207
$ws.Range("O2").Value = 3.341908864560384

# Sheet "2045" (sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("O2").Value = 5.381349308255924
